# Add a new "images" worksheet (map & progress bar tracking sheet) after the
# existing sheets, populate it with the image/update_date table, and make it
# the active sheet (matching the tabSelected / activeTab move in the diff).

$wb = $excel.ActiveWorkbook

# Grab a cell that already carries the "dd/mm/yyyy" style (s=3 in styles.xml)
# so the new date cells reuse that existing cell style instead of Excel
# fabricating a brand-new one.
$issueLog = $wb.Worksheets.Item("Issue Log")
$dateFormat = $issueLog.Range("E2").NumberFormat

# Create the new sheet as the last tab in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$imagesSheet = $wb.Worksheets.Add($null, $lastSheet)
$imagesSheet.Name = "images"

# Header row.
$imagesSheet.Range("A1").Value = "image"
$imagesSheet.Range("B1").Value = "update_date"
$imagesSheet.Range("A1:B1").Font.Bold = $true

# Data rows: image path + the date it was last updated (10-Apr-2025).
$images = @(
    "s05/plan.jpg",
    "s05/HA1-HA1'.jpg",
    "s05/HA2-HA2'.jpg",
    "s05/HA3-HA3'.jpg",
    "s05/HA4-HA4'.jpg",
    "s05/C12-C12'.jpg",
    "s05/C13-C13'.jpg",
    "s05/C13A-C13A'.jpg",
    "s05/C13B-C13B'.jpg"
)

$row = 2
foreach ($imagePath in $images) {
    $imagesSheet.Cells.Item($row, 1).Value = $imagePath
    $imagesSheet.Cells.Item($row, 2).Value = "4/10/2025"
    $row++
}

$lastRow = $row - 1
$imagesSheet.Range("B2:B" + $lastRow).NumberFormat = $dateFormat

# Size the columns to fit their contents.
$imagesSheet.Columns.Item(1).AutoFit()
$imagesSheet.Columns.Item(2).AutoFit()

# Match the portrait page setup used by the rest of the workbook's sheets.
$imagesSheet.PageSetup.Orientation = 1

# Put the selection/cursor on B1, and make "images" the active sheet/tab,
# mirroring the tabSelected move away from "Corridor Work".
[void]$imagesSheet.Range("B1").Select()
[void]$imagesSheet.Activate()
